# Update "想去人数" (number interested) values across sheets,
# matching the regenerated data snapshot at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibition) sheet
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 1340
$wsExhibition.Range("F6").Value = 370
$wsExhibition.Range("F7").Value = 3942
$wsExhibition.Range("F9").Value = 793
$wsExhibition.Range("F17").Value = 3456

# 演出 (Performance) sheet
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F12").Value = 8

# 本地生活 (Local Life) sheet
$wsLocalLife = $wb.Worksheets.Item("本地生活")
$wsLocalLife.Range("F4").Value = 2134

# 全部类型 (All Types) sheet - aggregated view of the above
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 2134
$wsAll.Range("F10").Value = 1340
$wsAll.Range("F16").Value = 370
$wsAll.Range("F17").Value = 3942
$wsAll.Range("F23").Value = 793
$wsAll.Range("F32").Value = 8
